# Apply the edits described by the diff:
# 1. Rename sheet "defect-report" to "Defect-report"
# 2. On "Test cases" sheet: remove scrolled topLeftCell (A56) -> scroll back to top (A1),
#    keep selection at E3
# 3. On "defect-report" sheet: change scroll position and selection
#    topLeftCell C3 -> A2 ; selection G8 -> C7

$wb = $excel.ActiveWorkbook

# --- Rename the "defect-report" sheet to "Defect-report" ---
$defectSheet = $wb.Worksheets.Item("defect-report")
$defectSheet.Name = "Defect-report"

# --- Fix up "Test cases" sheet view: scroll back to top (A1), keep selection E3 ---
$testCases = $wb.Worksheets.Item("Test cases")
$testCases.Activate() | Out-Null
$testCases.Range("E3").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- Fix up "Defect-report" sheet view: scroll top-left to A2, select C7 ---
$defectSheet.Activate() | Out-Null
$defectSheet.Range("C7").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
